# Update crypto price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.464.49"
$ws.Range("E2").Value = "  +0.84%  "

$ws.Range("D3").Value = "3.144.96"
$ws.Range("E3").Value = "  +0.40%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "603.24"
$ws.Range("E5").Value = "  -0.56%  "

$ws.Range("D6").Value = "144.09"
$ws.Range("E6").Value = "  -0.54%  "

$ws.Range("E7").Value = "  -0.26%  "

$ws.Range("D8").Value = "3.138.61"
$ws.Range("E8").Value = "  +0.64%  "

$ws.Range("E9").Value = "  +1.07%  "

$ws.Range("E10").Value = "  +0.29%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "5.40"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +3.35%  "

$ws.Range("D12").Value = "0.472"
$ws.Range("E12").Value = "  +0.47%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000256"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +2.30%  "

$ws.Range("D14").Value = "35.31"
$ws.Range("E14").Value = "  +0.54%  "

$ws.Range("D15").Value = "3.659.77"
$ws.Range("E15").Value = "  +0.43%  "

$ws.Range("E16").Value = "  +2.68%  "

$ws.Range("D17").Value = "64.388.89"
$ws.Range("E17").Value = "  +0.70%  "

$ws.Range("D18").Value = "3.135.99"
$ws.Range("E18").Value = "  +0.37%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.90"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.67%  "

$ws.Range("D20").Value = "483.21"
$ws.Range("E20").Value = "  +1.30%  "

$ws.Range("D21").Value = "14.62"
$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("D22").Value = "0.713"
$ws.Range("E22").Value = "  +0.93%  "

$ws.Range("D23").Value = "7.71"
$ws.Range("E23").Value = "  +0.27%  "

$ws.Range("D24").Value = "86.91"
$ws.Range("E24").Value = "  +3.38%  "

$ws.Range("D25").Value = "13.46"
$ws.Range("E25").Value = "  -0.41%  "

$ws.Range("D27").Value = "2.77"
$ws.Range("E27").Value = "  -0.42%  "

$ws.Range("D28").Value = "8.38"
$ws.Range("E28").Value = "  -0.47%  "

$ws.Range("D29").Value = "7.25"
$ws.Range("E29").Value = "  +6.57%  "

$ws.Range("E30").Value = "  -1.86%  "

$ws.Range("D31").Value = "0.113"
$ws.Range("E31").Value = "  -0.01%  "

$ws.Range("E32").Value = "  -0.11%  "

$ws.Range("D33").Value = "26.88"
$ws.Range("E33").Value = "  +2.52%  "

$ws.Range("D34").Value = "2.69"
$ws.Range("E34").Value = "  -0.30%  "

$ws.Range("E35").Value = "  -1.69%  "

$ws.Range("D36").Value = "6.02"
$ws.Range("E36").Value = "  +1.77%  "

$ws.Range("D37").Value = "0.0₃0763"
$ws.Range("E37").Value = "  +2.26%  "

$ws.Range("D38").Value = "52.64"
$ws.Range("E38").Value = "  -0.49%  "

$ws.Range("E39").Value = "  +2.89%  "

$ws.Range("D40").Value = "446.81"
$ws.Range("E40").Value = "  -2.39%  "

$ws.Range("D41").Value = "0.0394"
$ws.Range("E41").Value = "  +0.81%  "

$ws.Range("E42").Value = "  +1.46%  "

$ws.Range("D43").Value = "8.26"
$ws.Range("E43").Value = "  -0.83%  "

$ws.Range("D44").Value = "2.877.01"
$ws.Range("E44").Value = "  +1.07%  "

$ws.Range("D45").Value = "0.261"
$ws.Range("E45").Value = "  -1.21%  "

$ws.Range("D46").Value = "2.24"
$ws.Range("E46").Value = "  -0.81%  "

$ws.Range("D47").Value = "2.46"
$ws.Range("E47").Value = "  +1.60%  "

$ws.Range("D48").Value = "0.999"
$ws.Range("E48").Value = "  -0.02%  "

$ws.Range("D49").Value = "26.17"
$ws.Range("E49").Value = "  +0.56%  "

$ws.Range("E50").Value = "  +0.78%  "

$ws.Range("D51").Value = "121.68"
$ws.Range("E51").Value = "  +2.70%  "
